$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.985.23"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").Value = "2.101.50"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.90%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "347.86"
$ws.Range("E5").Value = "  +2.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5154"
$ws.Range("E7").Value = "  -1.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4421"
$ws.Range("E8").Value = "  -3.00%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09424"
$ws.Range("E9").Value = "  +3.28%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.44"
$ws.Range("E10").Value = "  -3.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.169"
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.30"
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("D13").Value = "2.098.09"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.734"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.183"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "98.90"
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001162"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.65"
$ws.Range("E19").Value = "  +6.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06662"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.229"
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").Value = "30.093.30"
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.63"
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.333"
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D26").Value = "2.352.59"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.95"
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.575"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.11"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.35"
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("E31").Value = "  -3.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1060"
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.646"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.223"
$ws.Range("E34").Value = "  -2.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.959"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.157"
$ws.Range("E36").Value = "  +4.89%  "
$ws.Range("E37").Value = "  -5.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02567"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06786"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2278"
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6885"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("E43").Value = "  +4.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6604"
$ws.Range("E44").Value = "  +1.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.21"
$ws.Range("E45").Value = "  -4.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.279"
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.632"
$ws.Range("E47").Value = "  -1.78%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000355"
$ws.Range("E48").Value = "  -4.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.219"
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.13"
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07219"
$ws.Range("E51").Value = "  -0.97%  "
